# "Finally excel download is fixed"
# - Rename sheet "A1"  -> "a1"
# - Rename sheet "a11" -> "A11"
# - Populate rows 2-12 of columns A:C on both of those sheets with the
#   (identical) data that the fixed download now produces.

$wb = $excel.ActiveWorkbook

$wsA1  = $wb.Worksheets.Item(2)   # currently named "A1"
$wsA11 = $wb.Worksheets.Item(3)   # currently named "a11"

$wsA1.Name  = "a1"
$wsA11.Name = "A11"

# Data for rows 2..12, columns A, B, C (B values must remain text)
$data = @(
    @("a1", "2", 1),
    @("A1", "2", 1),
    @("A1", "2", 1),
    @("A1", "2", 1),
    @("A1", "2", 1),
    @("A1", "2", 1),
    @("A1", "23", 1),
    @("A1", "23", 1),
    @("A1", "23", 1),
    @("a1", "21", 1),
    @("a1", "21", 1)
)

foreach ($ws in @($wsA1, $wsA11)) {
    # Keep column B as text so values like "2"/"23"/"21" aren't coerced to numbers
    $ws.Range("B2:B12").NumberFormat = "@"

    $r = 2
    foreach ($row in $data) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $r++
    }
}
